$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update CI(amplitude) text in L2
$ws.Range("L2").Value = "[0.3144185901385439, 0.4085361022475896]"

# Update p(amplitude) and q(amplitude) numeric values in M2 and N2
$ws.Range("M2").Value = [double]"3.177458296477198e-13"
$ws.Range("N2").Value = [double]"3.177458296477198e-13"

# Update CI(mesor) text in T2
$ws.Range("T2").Value = "[0.46740933878405044, 0.5219255848282057]"
